$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell that previously held "Good Morning" to the new text
$ws.Range("E8").Value = "GIT UPDATE"

# Update the selection on the sheet to reflect the active cell at save time
$ws.Range("E8").Select()
